# Update "Đơn 1 bác sĩ" sheet: insert a "Nhóm dịch vụ" column after F,
# and append 4 new columns (Tỉ lệ chiết khấu sale chính/phụ, Chiết khấu sale chính/phụ)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Đơn 1 bác sĩ")

# Insert a new column at G, shifting existing G..V to H..W
$ws1.Columns("G:G").Insert()

# New column G header + values
$ws1.Range("G1").Value = "Nhóm dịch vụ"
$ws1.Range("G2").Value = "Mũi"
# G3 (totals row) has no group value; leave blank like the other blank
# placeholder cells in that row (C3:H3, Q3:T3, etc.)

# New trailing columns X..AA (header row)
$ws1.Range("X1").Value = "Tỉ lệ chiết khấu sale chính"
$ws1.Range("Y1").Value = "Tỉ lệ chiết khấu sale phụ"
$ws1.Range("Z1").Value = "Chiết khấu sale chính"
$ws1.Range("AA1").Value = "Chiết khấu sale phụ"

# Row 2 values
$ws1.Range("X2").Value = 0.25
$ws1.Range("Y2").Value = 0
$ws1.Range("Z2").Value = 3682500
$ws1.Range("AA2").Value = 0

# Row 3 (totals) values
$ws1.Range("X3").Value = 0.25
$ws1.Range("Y3").Value = 0
$ws1.Range("Z3").Value = 3682500
$ws1.Range("AA3").Value = 0

# --- "Lương" sheet: populate salary / commission lookup parameters ---
$ws3 = $wb.Worksheets.Item("Lương")

$ws3.Range("A1").Value = "Danh mục"
$ws3.Range("B1").Value = 14

$ws3.Range("A2").Value = "Ngày công"
$ws3.Range("B2").Value = 11

$ws3.Range("A3").Value = "Phụ cấp"
$ws3.Range("B3").Value = 385000

$ws3.Range("A4").Value = "Lương cơ bản tại CẦN THƠ"
$ws3.Range("B4").Value = 0

$ws3.Range("A5").Value = "Chiết khấu sale chính tại CẦN THƠ"
$ws3.Range("B5").Value = 0

$ws3.Range("A6").Value = "Chiết khấu sale phụ tại CẦN THƠ"
$ws3.Range("B6").Value = 0

$ws3.Range("A7").Value = "Đơn 1 bác sĩ tại CẦN THƠ"
$ws3.Range("B7").Value = 1473000

$ws3.Range("A8").Value = "Đơn 2 bác sĩ tại CẦN THƠ"
$ws3.Range("B8").Value = 0

$ws3.Range("A9").Value = "Công phụ phẫu 1 tại CẦN THƠ"
$ws3.Range("B9").Value = 0

$ws3.Range("A10").Value = "Công phụ phẫu 2 tại CẦN THƠ"
$ws3.Range("B10").Value = 0

$ws3.Range("A11").Value = "Lương cơ bản tại LONG XUYÊN"
$ws3.Range("B11").Value = 0

$ws3.Range("A12").Value = "Chiết khấu sale chính tại LONG XUYÊN"
$ws3.Range("B12").Value = 0

$ws3.Range("A13").Value = "Chiết khấu sale phụ tại LONG XUYÊN"
$ws3.Range("B13").Value = 0

$ws3.Range("A14").Value = "Đơn 1 bác sĩ tại LONG XUYÊN"
$ws3.Range("B14").Value = 0

$ws3.Range("A15").Value = "Đơn 2 bác sĩ tại LONG XUYÊN"
$ws3.Range("B15").Value = 0

$ws3.Range("A16").Value = "Công phụ phẫu 1 tại LONG XUYÊN"
$ws3.Range("B16").Value = 0

$ws3.Range("A17").Value = "Công phụ phẫu 2 tại LONG XUYÊN"
$ws3.Range("B17").Value = 0

$ws3.Range("A18").Value = "Lương cơ bản tại SÓC TRĂNG"
$ws3.Range("B18").Value = 0

$ws3.Range("A19").Value = "Chiết khấu sale chính tại SÓC TRĂNG"
$ws3.Range("B19").Value = 0

$ws3.Range("A20").Value = "Chiết khấu sale phụ tại SÓC TRĂNG"
$ws3.Range("B20").Value = 0

$ws3.Range("A21").Value = "Đơn 1 bác sĩ tại SÓC TRĂNG"
$ws3.Range("B21").Value = 0

$ws3.Range("A22").Value = "Đơn 2 bác sĩ tại SÓC TRĂNG"
$ws3.Range("B22").Value = 0

$ws3.Range("A23").Value = "Công phụ phẫu 1 tại SÓC TRĂNG"
$ws3.Range("B23").Value = 0

$ws3.Range("A24").Value = "Công phụ phẫu 2 tại SÓC TRĂNG"
$ws3.Range("B24").Value = 0
